$d = $word.ActiveDocument

$d.Content.Find.Execute("34×26=884", $true, $false, $false, $false, $false, $true, 1, $false, "81×82=6642", 2) | Out-Null
$d.Content.Find.Execute("69×74=5106", $true, $false, $false, $false, $false, $true, 1, $false, "54×37=1998", 2) | Out-Null
$d.Content.Find.Execute("80×35=2800", $true, $false, $false, $false, $false, $true, 1, $false, "29×36=1044", 2) | Out-Null
$d.Content.Find.Execute("27×55=1485", $true, $false, $false, $false, $false, $true, 1, $false, "14×27=378", 2) | Out-Null
$d.Content.Find.Execute("47×81=3807", $true, $false, $false, $false, $false, $true, 1, $false, "13×54=702", 2) | Out-Null
$d.Content.Find.Execute("42×95=3990", $true, $false, $false, $false, $false, $true, 1, $false, "48×62=2976", 2) | Out-Null
$d.Content.Find.Execute("89×41=3649", $true, $false, $false, $false, $false, $true, 1, $false, "39×71=2769", 2) | Out-Null
$d.Content.Find.Execute("98×57=5586", $true, $false, $false, $false, $false, $true, 1, $false, "99×89=8811", 2) | Out-Null
$d.Content.Find.Execute("70×89=6230", $true, $false, $false, $false, $false, $true, 1, $false, "44×12=528", 2) | Out-Null
$d.Content.Find.Execute("53×100=5300", $true, $false, $false, $false, $false, $true, 1, $false, "14×40=560", 2) | Out-Null
$d.Content.Find.Execute("41×42=1722", $true, $false, $false, $false, $false, $true, 1, $false, "22×35=770", 2) | Out-Null
$d.Content.Find.Execute("51×70=3570", $true, $false, $false, $false, $false, $true, 1, $false, "14×21=294", 2) | Out-Null
$d.Content.Find.Execute("70×15=1050", $true, $false, $false, $false, $false, $true, 1, $false, "64×48=3072", 2) | Out-Null
$d.Content.Find.Execute("91×25=2275", $true, $false, $false, $false, $false, $true, 1, $false, "85×99=8415", 2) | Out-Null
$d.Content.Find.Execute("35×33=1155", $true, $false, $false, $false, $false, $true, 1, $false, "82×14=1148", 2) | Out-Null
$d.Content.Find.Execute("75×84=6300", $true, $false, $false, $false, $false, $true, 1, $false, "99×64=6336", 2) | Out-Null
$d.Content.Find.Execute("32×59=1888", $true, $false, $false, $false, $false, $true, 1, $false, "17×40=680", 2) | Out-Null
$d.Content.Find.Execute("25×55=1375", $true, $false, $false, $false, $false, $true, 1, $false, "41×98=4018", 2) | Out-Null
$d.Content.Find.Execute("99×69=6831", $true, $false, $false, $false, $false, $true, 1, $false, "63×48=3024", 2) | Out-Null
$d.Content.Find.Execute("80×76=6080", $true, $false, $false, $false, $false, $true, 1, $false, "98×27=2646", 2) | Out-Null
$d.Content.Find.Execute("84×72=6048", $true, $false, $false, $false, $false, $true, 1, $false, "76×52=3952", 2) | Out-Null
$d.Content.Find.Execute("52×35=1820", $true, $false, $false, $false, $false, $true, 1, $false, "46×44=2024", 2) | Out-Null
$d.Content.Find.Execute("66×24=1584", $true, $false, $false, $false, $false, $true, 1, $false, "53×32=1696", 2) | Out-Null
$d.Content.Find.Execute("44×25=1100", $true, $false, $false, $false, $false, $true, 1, $false, "89×18=1602", 2) | Out-Null
$d.Content.Find.Execute("62×70=4340", $true, $false, $false, $false, $false, $true, 1, $false, "74×49=3626", 2) | Out-Null
$d.Content.Find.Execute("10×41=410", $true, $false, $false, $false, $false, $true, 1, $false, "28×84=2352", 2) | Out-Null
$d.Content.Find.Execute("37×65=2405", $true, $false, $false, $false, $false, $true, 1, $false, "29×94=2726", 2) | Out-Null
$d.Content.Find.Execute("48×11=528", $true, $false, $false, $false, $false, $true, 1, $false, "25×16=400", 2) | Out-Null
$d.Content.Find.Execute("84×25=2100", $true, $false, $false, $false, $false, $true, 1, $false, "27×98=2646", 2) | Out-Null
$d.Content.Find.Execute("75×78=5850", $true, $false, $false, $false, $false, $true, 1, $false, "69×81=5589", 2) | Out-Null
$d.Content.Find.Execute("86×92=7912", $true, $false, $false, $false, $false, $true, 1, $false, "48×84=4032", 2) | Out-Null
$d.Content.Find.Execute("57×81=4617", $true, $false, $false, $false, $false, $true, 1, $false, "42×69=2898", 2) | Out-Null
$d.Content.Find.Execute("93×34=3162", $true, $false, $false, $false, $false, $true, 1, $false, "74×38=2812", 2) | Out-Null
$d.Content.Find.Execute("27×42=1134", $true, $false, $false, $false, $false, $true, 1, $false, "60×78=4680", 2) | Out-Null
$d.Content.Find.Execute("97×60=5820", $true, $false, $false, $false, $false, $true, 1, $false, "21×27=567", 2) | Out-Null
$d.Content.Find.Execute("43×44=1892", $true, $false, $false, $false, $false, $true, 1, $false, "82×67=5494", 2) | Out-Null
$d.Content.Find.Execute("80×10=800", $true, $false, $false, $false, $false, $true, 1, $false, "25×58=1450", 2) | Out-Null
$d.Content.Find.Execute("47×12=564", $true, $false, $false, $false, $false, $true, 1, $false, "49×99=4851", 2) | Out-Null
$d.Content.Find.Execute("59×46=2714", $true, $false, $false, $false, $false, $true, 1, $false, "33×86=2838", 2) | Out-Null
$d.Content.Find.Execute("69×54=3726", $true, $false, $false, $false, $false, $true, 1, $false, "72×98=7056", 2) | Out-Null
$d.Content.Find.Execute("91×23=2093", $true, $false, $false, $false, $false, $true, 1, $false, "89×66=5874", 2) | Out-Null
$d.Content.Find.Execute("68×34=2312", $true, $false, $false, $false, $false, $true, 1, $false, "25×33=825", 2) | Out-Null
$d.Content.Find.Execute("90×46=4140", $true, $false, $false, $false, $false, $true, 1, $false, "86×61=5246", 2) | Out-Null
$d.Content.Find.Execute("24×89=2136", $true, $false, $false, $false, $false, $true, 1, $false, "22×38=836", 2) | Out-Null
$d.Content.Find.Execute("51×92=4692", $true, $false, $false, $false, $false, $true, 1, $false, "98×97=9506", 2) | Out-Null
$d.Content.Find.Execute("70×100=7000", $true, $false, $false, $false, $false, $true, 1, $false, "44×82=3608", 2) | Out-Null
$d.Content.Find.Execute("21×89=1869", $true, $false, $false, $false, $false, $true, 1, $false, "78×44=3432", 2) | Out-Null
$d.Content.Find.Execute("83×70=5810", $true, $false, $false, $false, $false, $true, 1, $false, "32×68=2176", 2) | Out-Null
$d.Content.Find.Execute("91×75=6825", $true, $false, $false, $false, $false, $true, 1, $false, "97×74=7178", 2) | Out-Null
$d.Content.Find.Execute("36×97=3492", $true, $false, $false, $false, $false, $true, 1, $false, "28×57=1596", 2) | Out-Null
$d.Content.Find.Execute("86×21=1806", $true, $false, $false, $false, $false, $true, 1, $false, "43×19=817", 2) | Out-Null
$d.Content.Find.Execute("98×100=9800", $true, $false, $false, $false, $false, $true, 1, $false, "79×31=2449", 2) | Out-Null
$d.Content.Find.Execute("86×84=7224", $true, $false, $false, $false, $false, $true, 1, $false, "15×16=240", 2) | Out-Null
$d.Content.Find.Execute("97×55=5335", $true, $false, $false, $false, $false, $true, 1, $false, "26×29=754", 2) | Out-Null
$d.Content.Find.Execute("82×88=7216", $true, $false, $false, $false, $false, $true, 1, $false, "43×80=3440", 2) | Out-Null
$d.Content.Find.Execute("33×60=1980", $true, $false, $false, $false, $false, $true, 1, $false, "81×74=5994", 2) | Out-Null
$d.Content.Find.Execute("94×21=1974", $true, $false, $false, $false, $false, $true, 1, $false, "81×86=6966", 2) | Out-Null
$d.Content.Find.Execute("60×89=5340", $true, $false, $false, $false, $false, $true, 1, $false, "80×58=4640", 2) | Out-Null
$d.Content.Find.Execute("74×39=2886", $true, $false, $false, $false, $false, $true, 1, $false, "54×85=4590", 2) | Out-Null
$d.Content.Find.Execute("34×40=1360", $true, $false, $false, $false, $false, $true, 1, $false, "29×68=1972", 2) | Out-Null
$d.Content.Find.Execute("39×64=2496", $true, $false, $false, $false, $false, $true, 1, $false, "85×57=4845", 2) | Out-Null
$d.Content.Find.Execute("21×20=420", $true, $false, $false, $false, $false, $true, 1, $false, "20×64=1280", 2) | Out-Null
$d.Content.Find.Execute("72×40=2880", $true, $false, $false, $false, $false, $true, 1, $false, "92×95=8740", 2) | Out-Null
$d.Content.Find.Execute("45×29=1305", $true, $false, $false, $false, $false, $true, 1, $false, "77×52=4004", 2) | Out-Null
$d.Content.Find.Execute("34×63=2142", $true, $false, $false, $false, $false, $true, 1, $false, "59×95=5605", 2) | Out-Null
$d.Content.Find.Execute("30×96=2880", $true, $false, $false, $false, $false, $true, 1, $false, "18×35=630", 2) | Out-Null
$d.Content.Find.Execute("45×85=3825", $true, $false, $false, $false, $false, $true, 1, $false, "20×69=1380", 2) | Out-Null
$d.Content.Find.Execute("33×36=1188", $true, $false, $false, $false, $false, $true, 1, $false, "51×81=4131", 2) | Out-Null
$d.Content.Find.Execute("53×92=4876", $true, $false, $false, $false, $false, $true, 1, $false, "18×54=972", 2) | Out-Null
$d.Content.Find.Execute("74×35=2590", $true, $false, $false, $false, $false, $true, 1, $false, "65×49=3185", 2) | Out-Null
$d.Content.Find.Execute("82×28=2296", $true, $false, $false, $false, $false, $true, 1, $false, "55×90=4950", 2) | Out-Null
$d.Content.Find.Execute("27×90=2430", $true, $false, $false, $false, $false, $true, 1, $false, "33×31=1023", 2) | Out-Null
$d.Content.Find.Execute("76×15=1140", $true, $false, $false, $false, $false, $true, 1, $false, "74×68=5032", 2) | Out-Null
$d.Content.Find.Execute("41×99=4059", $true, $false, $false, $false, $false, $true, 1, $false, "16×41=656", 2) | Out-Null
$d.Content.Find.Execute("80×66=5280", $true, $false, $false, $false, $false, $true, 1, $false, "24×21=504", 2) | Out-Null
$d.Content.Find.Execute("78×61=4758", $true, $false, $false, $false, $false, $true, 1, $false, "15×78=1170", 2) | Out-Null
$d.Content.Find.Execute("16×46=736", $true, $false, $false, $false, $false, $true, 1, $false, "62×29=1798", 2) | Out-Null
$d.Content.Find.Execute("63×66=4158", $true, $false, $false, $false, $false, $true, 1, $false, "47×56=2632", 2) | Out-Null
$d.Content.Find.Execute("81×60=4860", $true, $false, $false, $false, $false, $true, 1, $false, "87×76=6612", 2) | Out-Null
$d.Content.Find.Execute("52×21=1092", $true, $false, $false, $false, $false, $true, 1, $false, "88×69=6072", 2) | Out-Null
$d.Content.Find.Execute("65×12=780", $true, $false, $false, $false, $false, $true, 1, $false, "12×20=240", 2) | Out-Null
$d.Content.Find.Execute("35×53=1855", $true, $false, $false, $false, $false, $true, 1, $false, "15×95=1425", 2) | Out-Null
$d.Content.Find.Execute("13×82=1066", $true, $false, $false, $false, $false, $true, 1, $false, "24×17=408", 2) | Out-Null
$d.Content.Find.Execute("71×84=5964", $true, $false, $false, $false, $false, $true, 1, $false, "44×58=2552", 2) | Out-Null
$d.Content.Find.Execute("35×68=2380", $true, $false, $false, $false, $false, $true, 1, $false, "91×37=3367", 2) | Out-Null
$d.Content.Find.Execute("23×62=1426", $true, $false, $false, $false, $false, $true, 1, $false, "17×56=952", 2) | Out-Null
$d.Content.Find.Execute("49×43=2107", $true, $false, $false, $false, $false, $true, 1, $false, "26×84=2184", 2) | Out-Null
$d.Content.Find.Execute("29×70=2030", $true, $false, $false, $false, $false, $true, 1, $false, "84×30=2520", 2) | Out-Null
$d.Content.Find.Execute("80×75=6000", $true, $false, $false, $false, $false, $true, 1, $false, "29×91=2639", 2) | Out-Null
$d.Content.Find.Execute("78×54=4212", $true, $false, $false, $false, $false, $true, 1, $false, "81×35=2835", 2) | Out-Null
$d.Content.Find.Execute("57×70=3990", $true, $false, $false, $false, $false, $true, 1, $false, "54×14=756", 2) | Out-Null
$d.Content.Find.Execute("65×69=4485", $true, $false, $false, $false, $false, $true, 1, $false, "38×38=1444", 2) | Out-Null
$d.Content.Find.Execute("61×70=4270", $true, $false, $false, $false, $false, $true, 1, $false, "16×33=528", 2) | Out-Null
$d.Content.Find.Execute("13×48=624", $true, $false, $false, $false, $false, $true, 1, $false, "60×57=3420", 2) | Out-Null
$d.Content.Find.Execute("12×70=840", $true, $false, $false, $false, $false, $true, 1, $false, "12×19=228", 2) | Out-Null
$d.Content.Find.Execute("69×85=5865", $true, $false, $false, $false, $false, $true, 1, $false, "43×55=2365", 2) | Out-Null
$d.Content.Find.Execute("100×11=1100", $true, $false, $false, $false, $false, $true, 1, $false, "66×40=2640", 2) | Out-Null
$d.Content.Find.Execute("95×94=8930", $true, $false, $false, $false, $false, $true, 1, $false, "55×99=5445", 2) | Out-Null
$d.Content.Find.Execute("73×17=1241", $true, $false, $false, $false, $false, $true, 1, $false, "26×70=1820", 2) | Out-Null
$d.Content.Find.Execute("55×57=3135", $true, $false, $false, $false, $false, $true, 1, $false, "24×30=720", 2) | Out-Null
